# Scheduled-runner style refresh of cached market/profit figures across the
# FFXIV Leve-profit sheets (columns H-N: average prices, leve NQ/HQ prices,
# and NQ/HQ profit). Values are plain numbers (no formulas in the source
# workbook), so each changed cell is written directly; cells that the
# refresh run dropped (no longer computable) are cleared so they disappear
# from the row entirely, and cells newly introduced are written fresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 893.5  # H12
$ws.Cells.Item(12, 9).Value = 476.57144  # I12
$ws.Cells.Item(12, 10).Value = 1866.3334  # J12
$ws.Cells.Item(12, 11).Value = 476.57144  # K12
$ws.Cells.Item(12, 12).Value = 1866.3334  # L12
$ws.Cells.Item(12, 13).Value = -306.57144  # M12
$ws.Cells.Item(12, 14).Value = -2206.3334  # N12
$ws.Cells.Item(31, 8).Value = 1000  # H31
$ws.Cells.Item(31, 9).Value = 0  # I31
$ws.Cells.Item(31, 11).Value = 0  # K31
$ws.Cells.Item(31, 13).ClearContents()  # M31 was 221
$ws.Cells.Item(55, 8).Value = 466.66666  # H55
$ws.Cells.Item(55, 9).Value = 500  # I55
$ws.Cells.Item(55, 10).Value = 400  # J55
$ws.Cells.Item(55, 11).Value = 500  # K55
$ws.Cells.Item(55, 12).Value = 400  # L55
$ws.Cells.Item(55, 13).Value = -286  # M55
$ws.Cells.Item(55, 14).Value = -828  # N55
$ws.Cells.Item(64, 8).Value = 9444.444  # H64
$ws.Cells.Item(64, 9).Value = 4533.3335  # I64
$ws.Cells.Item(64, 11).Value = 4533.3335  # K64
$ws.Cells.Item(64, 13).Value = -4285.3335  # M64
$ws.Cells.Item(67, 8).Value = 9444.444  # H67
$ws.Cells.Item(67, 9).Value = 4533.3335  # I67
$ws.Cells.Item(67, 11).Value = 4533.3335  # K67
$ws.Cells.Item(67, 13).Value = -3675.3335  # M67
$ws.Cells.Item(92, 8).Value = 548.8889  # H92
$ws.Cells.Item(92, 9).Value = 491.875  # I92
$ws.Cells.Item(92, 10).Value = 1005  # J92
$ws.Cells.Item(92, 11).Value = 491.875  # K92
$ws.Cells.Item(92, 12).Value = 1005  # L92
$ws.Cells.Item(92, 13).Value = 756.125  # M92
$ws.Cells.Item(92, 14).Value = -3501  # N92
$ws.Cells.Item(132, 8).Value = 1190.3334  # H132
$ws.Cells.Item(132, 9).Value = 1189.9131  # I132
$ws.Cells.Item(132, 10).Value = 1200  # J132
$ws.Cells.Item(132, 11).Value = 3569.7393  # K132
$ws.Cells.Item(132, 12).Value = 3600  # L132
$ws.Cells.Item(132, 13).Value = -1039.7393  # M132
$ws.Cells.Item(132, 14).Value = -8660  # N132
$ws.Cells.Item(137, 8).Value = 1850.3704  # H137
$ws.Cells.Item(137, 9).Value = 1657.3182  # I137
$ws.Cells.Item(137, 11).Value = 4971.9546  # K137
$ws.Cells.Item(137, 13).Value = -2421.9546  # M137
$ws.Cells.Item(138, 8).Value = 28574670  # H138
$ws.Cells.Item(138, 9).Value = 2912.7273  # I138
$ws.Cells.Item(138, 10).Value = 41670060  # J138
$ws.Cells.Item(138, 11).Value = 8738.1819  # K138
$ws.Cells.Item(138, 12).Value = 125010180  # L138
$ws.Cells.Item(138, 13).Value = -3598.1819  # M138
$ws.Cells.Item(138, 14).Value = -125020460  # N138

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(13, 8).Value = 0  # H13
$ws.Cells.Item(13, 10).Value = 0  # J13
$ws.Cells.Item(13, 12).Value = 0  # L13
$ws.Cells.Item(13, 14).ClearContents()  # N13 was -2788
$ws.Cells.Item(32, 8).Value = 6114.3413  # H32
$ws.Cells.Item(32, 9).Value = 6114.3413  # I32
$ws.Cells.Item(32, 11).Value = 6114.3413  # K32
$ws.Cells.Item(32, 13).Value = -5827.3413  # M32
$ws.Cells.Item(45, 8).Value = 8091.769  # H45
$ws.Cells.Item(45, 9).Value = 3799.75  # I45
$ws.Cells.Item(45, 10).Value = 9999.333000000001  # J45
$ws.Cells.Item(45, 11).Value = 3799.75  # K45
$ws.Cells.Item(45, 12).Value = 9999.333000000001  # L45
$ws.Cells.Item(45, 13).Value = -3422.75  # M45
$ws.Cells.Item(45, 14).Value = -10753.333  # N45
$ws.Cells.Item(61, 8).Value = 26322594  # H61
$ws.Cells.Item(61, 9).Value = 33338420  # I61
$ws.Cells.Item(61, 11).Value = 33338420  # K61
$ws.Cells.Item(61, 13).Value = -33338208  # M61
$ws.Cells.Item(97, 8).Value = 494.9565  # H97
$ws.Cells.Item(97, 9).Value = 587.5333000000001  # I97
$ws.Cells.Item(97, 11).Value = 587.5333000000001  # K97
$ws.Cells.Item(97, 13).Value = -91.53330000000005  # M97
$ws.Cells.Item(131, 8).Value = 39331.668  # H131
$ws.Cells.Item(131, 10).Value = 39331.668  # J131
$ws.Cells.Item(131, 12).Value = 39331.668  # L131
$ws.Cells.Item(131, 14).Value = -49411.668  # N131
$ws.Cells.Item(132, 8).Value = 4144.93  # H132
$ws.Cells.Item(132, 9).Value = 3748.3057  # I132
$ws.Cells.Item(132, 10).Value = 6184.7144  # J132
$ws.Cells.Item(132, 11).Value = 11244.9171  # K132
$ws.Cells.Item(132, 12).Value = 18554.1432  # L132
$ws.Cells.Item(132, 13).Value = -8714.917099999999  # M132
$ws.Cells.Item(132, 14).Value = -23614.1432  # N132
$ws.Cells.Item(136, 8).Value = 26322594  # H136
$ws.Cells.Item(136, 9).Value = 33338420  # I136
$ws.Cells.Item(136, 11).Value = 100015260  # K136
$ws.Cells.Item(136, 13).Value = -100012710  # M136

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5, 8).Value = 2999.5  # H5
$ws.Cells.Item(5, 9).Value = 2999.5  # I5
$ws.Cells.Item(5, 10).Value = 0  # J5
$ws.Cells.Item(5, 11).Value = 2999.5  # K5
$ws.Cells.Item(5, 12).Value = 0  # L5
$ws.Cells.Item(5, 13).Value = -2886.5  # M5
$ws.Cells.Item(5, 14).ClearContents()  # N5 was -4226
$ws.Cells.Item(22, 8).Value = 1980  # H22
$ws.Cells.Item(22, 9).Value = 1980  # I22
$ws.Cells.Item(22, 11).Value = 1980  # K22
$ws.Cells.Item(22, 13).Value = -1807  # M22
$ws.Cells.Item(36, 8).Value = 4317  # H36
$ws.Cells.Item(36, 9).Value = 4317  # I36
$ws.Cells.Item(36, 11).Value = 4317  # K36
$ws.Cells.Item(36, 13).Value = -3783  # M36

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2486  # H132
$ws.Cells.Item(132, 9).Value = 1820.2667  # I132
$ws.Cells.Item(132, 11).Value = 5460.800099999999  # K132
$ws.Cells.Item(132, 13).Value = -2930.800099999999  # M132

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 5079  # H3
$ws.Cells.Item(3, 9).Value = 5288.154  # I3
$ws.Cells.Item(3, 10).Value = 4399.25  # J3
$ws.Cells.Item(3, 11).Value = 15864.462  # K3
$ws.Cells.Item(3, 12).Value = 13197.75  # L3
$ws.Cells.Item(3, 13).Value = -15752.462  # M3
$ws.Cells.Item(3, 14).Value = -13421.75  # N3
$ws.Cells.Item(39, 8).Value = 1375  # H39
$ws.Cells.Item(39, 10).Value = 1428.5714  # J39
$ws.Cells.Item(39, 12).Value = 4285.7142  # L39
$ws.Cells.Item(39, 14).Value = -4873.7142  # N39
$ws.Cells.Item(55, 8).Value = 823.5  # H55
$ws.Cells.Item(55, 10).Value = 999.875  # J55
$ws.Cells.Item(55, 12).Value = 2999.625  # L55
$ws.Cells.Item(55, 14).Value = -3353.625  # N55
$ws.Cells.Item(87, 8).Value = 4522.3335  # H87
$ws.Cells.Item(87, 9).Value = 4522.3335  # I87
$ws.Cells.Item(87, 11).Value = 13567.0005  # K87
$ws.Cells.Item(87, 13).Value = -12319.0005  # M87
$ws.Cells.Item(90, 8).Value = 4522.3335  # H90
$ws.Cells.Item(90, 9).Value = 4522.3335  # I90
$ws.Cells.Item(90, 11).Value = 40701.0015  # K90
$ws.Cells.Item(90, 13).Value = -34461.0015  # M90
$ws.Cells.Item(107, 8).Value = 1483.1666  # H107
$ws.Cells.Item(107, 9).Value = 799.6667  # I107
$ws.Cells.Item(107, 11).Value = 2399.0001  # K107
$ws.Cells.Item(107, 13).Value = -479.0001000000002  # M107
$ws.Cells.Item(117, 8).Value = 1001668.9  # H117
$ws.Cells.Item(117, 10).Value = 1430347.2  # J117
$ws.Cells.Item(117, 12).Value = 4291041.6  # L117
$ws.Cells.Item(117, 14).Value = -4297925.6  # N117
$ws.Cells.Item(139, 8).Value = 2159.5625  # H139
$ws.Cells.Item(139, 9).Value = 1853.5667  # I139
$ws.Cells.Item(139, 11).Value = 5560.7001  # K139
$ws.Cells.Item(139, 13).Value = -420.7001  # M139

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 30054  # H15
$ws.Cells.Item(15, 10).Value = 30054  # J15
$ws.Cells.Item(15, 12).Value = 30054  # L15
$ws.Cells.Item(15, 14).Value = -30630  # N15
$ws.Cells.Item(41, 8).Value = 11697.167  # H41
$ws.Cells.Item(41, 9).Value = 12306.091  # I41
$ws.Cells.Item(41, 10).Value = 4999  # J41
$ws.Cells.Item(41, 11).Value = 12306.091  # K41
$ws.Cells.Item(41, 12).Value = 4999  # L41
$ws.Cells.Item(41, 13).Value = -11951.091  # M41
$ws.Cells.Item(41, 14).Value = -5709  # N41
$ws.Cells.Item(80, 8).Value = 3794.5  # H80
$ws.Cells.Item(80, 9).Value = 3305.9333  # I80
$ws.Cells.Item(80, 10).Value = 4608.778  # J80
$ws.Cells.Item(80, 11).Value = 3305.9333  # K80
$ws.Cells.Item(80, 12).Value = 4608.778  # L80
$ws.Cells.Item(80, 13).Value = -2307.9333  # M80
$ws.Cells.Item(80, 14).Value = -6604.778  # N80
$ws.Cells.Item(81, 8).Value = 30054  # H81
$ws.Cells.Item(81, 10).Value = 30054  # J81
$ws.Cells.Item(81, 12).Value = 30054  # L81
$ws.Cells.Item(81, 14).Value = -32050  # N81
$ws.Cells.Item(83, 8).Value = 3794.5  # H83
$ws.Cells.Item(83, 9).Value = 3305.9333  # I83
$ws.Cells.Item(83, 10).Value = 4608.778  # J83
$ws.Cells.Item(83, 11).Value = 16529.6665  # K83
$ws.Cells.Item(83, 12).Value = 23043.89  # L83
$ws.Cells.Item(83, 13).Value = -11537.6665  # M83
$ws.Cells.Item(83, 14).Value = -33027.89  # N83
$ws.Cells.Item(84, 8).Value = 30054  # H84
$ws.Cells.Item(84, 10).Value = 30054  # J84
$ws.Cells.Item(84, 12).Value = 90162  # L84
$ws.Cells.Item(84, 14).Value = -100146  # N84
$ws.Cells.Item(102, 8).Value = 8143.077  # H102
$ws.Cells.Item(102, 9).Value = 6896.8184  # I102
$ws.Cells.Item(102, 10).Value = 14997.5  # J102
$ws.Cells.Item(102, 11).Value = 6896.8184  # K102
$ws.Cells.Item(102, 12).Value = 14997.5  # L102
$ws.Cells.Item(102, 13).Value = -5274.8184  # M102
$ws.Cells.Item(102, 14).Value = -18241.5  # N102
$ws.Cells.Item(107, 8).Value = 400.27777  # H107
$ws.Cells.Item(107, 9).Value = 423.36365  # I107
$ws.Cells.Item(107, 10).Value = 364  # J107
$ws.Cells.Item(107, 11).Value = 423.36365  # K107
$ws.Cells.Item(107, 12).Value = 364  # L107
$ws.Cells.Item(107, 13).Value = 1496.63635  # M107
$ws.Cells.Item(107, 14).Value = -4204  # N107
$ws.Cells.Item(126, 8).Value = 2775.4348  # H126
$ws.Cells.Item(126, 9).Value = 1764.4615  # I126
$ws.Cells.Item(126, 11).Value = 5293.3845  # K126
$ws.Cells.Item(126, 13).Value = -2823.3845  # M126
$ws.Cells.Item(134, 8).Value = 61996.332  # H134
$ws.Cells.Item(134, 10).Value = 61996.332  # J134
$ws.Cells.Item(134, 12).Value = 185988.996  # L134
$ws.Cells.Item(134, 14).Value = -191058.996  # N134

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(75, 8).Value = 0  # H75
$ws.Cells.Item(75, 10).Value = 0  # J75
$ws.Cells.Item(75, 12).Value = 0  # L75
$ws.Cells.Item(75, 14).ClearContents()  # N75 was -16547.5
$ws.Cells.Item(78, 8).Value = 0  # H78
$ws.Cells.Item(78, 10).Value = 0  # J78
$ws.Cells.Item(78, 12).Value = 0  # L78
$ws.Cells.Item(78, 14).ClearContents()  # N78 was -53386.5
$ws.Cells.Item(122, 8).Value = 2394  # H122
$ws.Cells.Item(122, 9).Value = 1876.25  # I122
$ws.Cells.Item(122, 11).Value = 5628.75  # K122
$ws.Cells.Item(122, 13).Value = -3178.75  # M122
$ws.Cells.Item(132, 8).Value = 9849.736999999999  # H132
$ws.Cells.Item(132, 9).Value = 10859.1875  # I132
$ws.Cells.Item(132, 11).Value = 32577.5625  # K132
$ws.Cells.Item(132, 13).Value = -30047.5625  # M132

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 4792.9  # H4
$ws.Cells.Item(4, 9).Value = 3918.5715  # I4
$ws.Cells.Item(4, 10).Value = 6833  # J4
$ws.Cells.Item(4, 11).Value = 3918.5715  # K4
$ws.Cells.Item(4, 12).Value = 6833  # L4
$ws.Cells.Item(4, 13).Value = -3805.5715  # M4
$ws.Cells.Item(4, 14).Value = -7059  # N4
$ws.Cells.Item(32, 8).Value = 2833.25  # H32
$ws.Cells.Item(32, 10).Value = 3666.5  # J32
$ws.Cells.Item(32, 12).Value = 3666.5  # L32
$ws.Cells.Item(32, 14).Value = -4300.5  # N32
$ws.Cells.Item(107, 8).Value = 786.4643  # H107
$ws.Cells.Item(107, 9).Value = 679.26086  # I107
$ws.Cells.Item(107, 10).Value = 1279.6  # J107
$ws.Cells.Item(107, 11).Value = 2037.78258  # K107
$ws.Cells.Item(107, 12).Value = 3838.8  # L107
$ws.Cells.Item(107, 13).Value = -117.7825800000001  # M107
$ws.Cells.Item(107, 14).Value = -7678.799999999999  # N107
$ws.Cells.Item(132, 8).Value = 2204.8857  # H132
$ws.Cells.Item(132, 9).Value = 1931.5  # I132
$ws.Cells.Item(132, 10).Value = 2994.6667  # J132
$ws.Cells.Item(132, 11).Value = 5794.5  # K132
$ws.Cells.Item(132, 12).Value = 8984.000100000001  # L132
$ws.Cells.Item(132, 13).Value = -3264.5  # M132
$ws.Cells.Item(132, 14).Value = -14044.0001  # N132
$ws.Cells.Item(136, 8).Value = 3944.4211  # H136
$ws.Cells.Item(136, 9).Value = 2424.4167  # I136
$ws.Cells.Item(136, 10).Value = 6550.143  # J136
$ws.Cells.Item(136, 11).Value = 7273.250100000001  # K136
$ws.Cells.Item(136, 12).Value = 19650.429  # L136
$ws.Cells.Item(136, 13).Value = -4723.250100000001  # M136
$ws.Cells.Item(136, 14).Value = -24750.429  # N136
